$d = $word.ActiveDocument

# Find the paragraph that ends with "...was unsuccessful." and insert a new
# paragraph right after it containing the "Ruled out..." sentence.
$found = $d.Content.Find.Execute(
    "Tried to duplicate in the dev database but was unsuccessful.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tried to duplicate in the dev database but was unsuccessful.^p" +
    "Ruled out that the error occurred as a result of selling an item that was stocked in another location.",
    2)
